$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the old "Programa resumido:" row (row 12),
# shifting the existing rows (old 12-20) down to 15-23.
$ws.Range("A12:A14").EntireRow.Insert()

# The row-insert carries the column formatting (styles 1/2/3) into every
# cell of the new rows; clear the ones that should stay genuinely empty so
# they are not serialized as blank styled cells.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# New row 12: section header "Docentes responsáveis:" (column A only).
$ws.Range("A12").Value = "Docentes responsáveis:"

# New row 13: first responsible professor (columns B and C, same text).
$ws.Range("B13").Value = "5817344 - Livia Melo Carneiro"
$ws.Range("C13").Value = "5817344 - Livia Melo Carneiro"

# New row 14: second responsible professor (columns B and C, same text).
$ws.Range("B14").Value = "6310296 - Patrícia Caroline Molgero Da Rós"
$ws.Range("C14").Value = "6310296 - Patrícia Caroline Molgero Da Rós"
